# Fruta / hortaliza, semanal
# Insert a new weekly price record above the existing row 269
# (shifting the subsequent rows down by one) on the "Ajo" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 269; everything from the old
# row 269 onward shifts down one row (old 269 -> 270, ..., old 286 -> 287).
$ws.Rows.Item(269).Insert()

# Populate the newly inserted row 269 with the new weekly record.
$ws.Range("A269").Value = 8
$ws.Range("B269").Value = "Terminal La Palmera de La Serena"
$ws.Range("C269").Value = "Coquimbo"
$ws.Range("D269").Value = 44746
$ws.Range("E269").Value = 4
$ws.Range("F269").Value = 100112003
$ws.Range("G269").Value = "Ajo"
$ws.Range("H269").Value = "Chino"
$ws.Range("I269").Value = "Primera"
$ws.Range("J269").Value = 480
$ws.Range("K269").Value = 19000
$ws.Range("L269").Value = 20000
$ws.Range("M269").Value = 19500
$ws.Range("N269").Value = "`$/caja 10 kilos"
$ws.Range("O269").Value = "China"
$ws.Range("P269").Value = 1950
$ws.Range("Q269").Value = 10
$ws.Range("R269").Value = "Hortaliza"
